$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1,1).Value2 = 'Datos actualizados a 19 de Julio de 2020 a las 11:50'

# Row 4
$ws.Cells.Item(4,2).Value2 = 3833716
$ws.Cells.Item(4,3).Value2 = 445
$ws.Cells.Item(4,4).Value2 = 1775450
$ws.Cells.Item(4,5).Value2 = 1915385
$ws.Cells.Item(4,7).Value2 = 4
$ws.Cells.Item(4,8).Value2 = 142881

# Row 6
$ws.Cells.Item(6,2).Value2 = 1078782
$ws.Cells.Item(6,3).Value2 = 918
$ws.Cells.Item(6,5).Value2 = 374088

# Row 19
$ws.Cells.Item(19,1).Value2 = 'Banglades'
$ws.Cells.Item(19,2).Value2 = 204525
$ws.Cells.Item(19,3).Value2 = 2459
$ws.Cells.Item(19,4).Value2 = 111642
$ws.Cells.Item(19,5).Value2 = 90265
$ws.Cells.Item(19,7).Value2 = 37
$ws.Cells.Item(19,8).Value2 = 2618

# Row 20
$ws.Cells.Item(20,1).Value2 = 'Alemania'
$ws.Cells.Item(20,2).Value2 = 202572
$ws.Cells.Item(20,4).Value2 = 187800
$ws.Cells.Item(20,5).Value2 = 5610
$ws.Cells.Item(20,8).Value2 = 9162

# Row 28
$ws.Cells.Item(28,2).Value2 = 86521
$ws.Cells.Item(28,3).Value2 = 1639
$ws.Cells.Item(28,4).Value2 = 45401
$ws.Cells.Item(28,5).Value2 = 36977
$ws.Cells.Item(28,7).Value2 = 127
$ws.Cells.Item(28,8).Value2 = 4143

# Row 34
$ws.Cells.Item(34,1).Value2 = 'Oman'
$ws.Cells.Item(34,2).Value2 = 66661
$ws.Cells.Item(34,3).Value2 = 1157
$ws.Cells.Item(34,4).Value2 = 44004
$ws.Cells.Item(34,5).Value2 = 22339
$ws.Cells.Item(34,7).Value2 = 10
$ws.Cells.Item(34,8).Value2 = 318

# Row 35
$ws.Cells.Item(35,1).Value2 = 'Bielorrusia'
$ws.Cells.Item(35,2).Value2 = 65953
$ws.Cells.Item(35,4).Value2 = 57856
$ws.Cells.Item(35,5).Value2 = 7602
$ws.Cells.Item(35,8).Value2 = 495

# Row 47
$ws.Cells.Item(47,2).Value2 = 40104
$ws.Cells.Item(47,3).Value2 = 358
$ws.Cells.Item(47,5).Value2 = 8188
$ws.Cells.Item(47,7).Value2 = 6
$ws.Cells.Item(47,8).Value2 = 1624

# Row 52
$ws.Cells.Item(52,2).Value2 = 35475
$ws.Cells.Item(52,3).Value2 = 174
$ws.Cells.Item(52,4).Value2 = 23634
$ws.Cells.Item(52,5).Value2 = 10660
$ws.Cells.Item(52,7).Value2 = 17
$ws.Cells.Item(52,8).Value2 = 1181

# Row 64
$ws.Cells.Item(64,2).Value2 = 19655
$ws.Cells.Item(64,3).Value2 = 82
$ws.Cells.Item(64,4).Value2 = 17599
$ws.Cells.Item(64,5).Value2 = 1345

# Row 74
$ws.Cells.Item(74,1).Value2 = 'El Salvador'
$ws.Cells.Item(74,2).Value2 = 11846
$ws.Cells.Item(74,3).Value2 = 338
$ws.Cells.Item(74,4).Value2 = 6705
$ws.Cells.Item(74,5).Value2 = 4806
$ws.Cells.Item(74,7).Value2 = 11
$ws.Cells.Item(74,8).Value2 = 335

# Row 75
$ws.Cells.Item(75,1).Value2 = 'Australia'
$ws.Cells.Item(75,2).Value2 = 11802
$ws.Cells.Item(75,3).Value2 = 361
$ws.Cells.Item(75,4).Value2 = 8273
$ws.Cells.Item(75,5).Value2 = 3407
$ws.Cells.Item(75,7).Value2 = 4
$ws.Cells.Item(75,8).Value2 = 122

# Row 82
$ws.Cells.Item(82,2).Value2 = 8779
$ws.Cells.Item(82,3).Value2 = 15
$ws.Cells.Item(82,4).Value2 = 8553
$ws.Cells.Item(82,5).Value2 = 103
$ws.Cells.Item(82,7).Value2 = 1
$ws.Cells.Item(82,8).Value2 = 123

# Row 85
$ws.Cells.Item(85,1).Value2 = 'Estado de Palestina'
$ws.Cells.Item(85,2).Value2 = 8549
$ws.Cells.Item(85,3).Value2 = 345
$ws.Cells.Item(85,4).Value2 = 1921
$ws.Cells.Item(85,5).Value2 = 6569
$ws.Cells.Item(85,8).Value2 = 59

# Row 86
$ws.Cells.Item(86,1).Value2 = 'Consejo Danes para los Refugiados'
$ws.Cells.Item(86,2).Value2 = 8403
$ws.Cells.Item(86,3).Value2 = 79
$ws.Cells.Item(86,4).Value2 = 4335
$ws.Cells.Item(86,5).Value2 = 3874
$ws.Cells.Item(86,7).Value2 = 1
$ws.Cells.Item(86,8).Value2 = 194

# Row 88
$ws.Cells.Item(88,2).Value2 = 7335
$ws.Cells.Item(88,3).Value2 = 17
$ws.Cells.Item(88,5).Value2 = 127

# Row 110
$ws.Cells.Item(110,2).Value2 = 2907
$ws.Cells.Item(110,3).Value2 = 97
$ws.Cells.Item(110,4).Value2 = 1135
$ws.Cells.Item(110,5).Value2 = 1713
$ws.Cells.Item(110,7).Value2 = 4
$ws.Cells.Item(110,8).Value2 = 59

# Row 113
$ws.Cells.Item(113,2).Value2 = 2708
$ws.Cells.Item(113,3).Value2 = 4
$ws.Cells.Item(113,4).Value2 = 2035
$ws.Cells.Item(113,5).Value2 = 662

# Row 123
$ws.Cells.Item(123,2).Value2 = 1946
$ws.Cells.Item(123,3).Value2 = 6
$ws.Cells.Item(123,5).Value2 = 266
$ws.Cells.Item(123,7).Value2 = 1
$ws.Cells.Item(123,8).Value2 = 112

# Row 124
$ws.Cells.Item(124,1).Value2 = 'Lituania'
$ws.Cells.Item(124,2).Value2 = 1932
$ws.Cells.Item(124,3).Value2 = 17
$ws.Cells.Item(124,4).Value2 = 1600
$ws.Cells.Item(124,5).Value2 = 252
$ws.Cells.Item(124,8).Value2 = 80

# Row 125
$ws.Cells.Item(125,1).Value2 = 'Islandia'
$ws.Cells.Item(125,2).Value2 = 1922
$ws.Cells.Item(125,4).Value2 = 1902
$ws.Cells.Item(125,5).Value2 = 10
$ws.Cells.Item(125,8).Value2 = 10

# Row 126
$ws.Cells.Item(126,1).Value2 = 'Hong Kong'
$ws.Cells.Item(126,2).Value2 = 1886
$ws.Cells.Item(126,3).Value2 = 108
$ws.Cells.Item(126,4).Value2 = 1294
$ws.Cells.Item(126,5).Value2 = 580
$ws.Cells.Item(126,8).Value2 = 12

# Row 127
$ws.Cells.Item(127,1).Value2 = 'Libia'
$ws.Cells.Item(127,2).Value2 = 1791
$ws.Cells.Item(127,4).Value2 = 385
$ws.Cells.Item(127,5).Value2 = 1358
$ws.Cells.Item(127,8).Value2 = 48

# Row 142
$ws.Cells.Item(142,2).Value2 = 1065
$ws.Cells.Item(142,3).Value2 = 3
$ws.Cells.Item(142,5).Value2 = 42

# Row 165
$ws.Cells.Item(165,2).Value2 = 341
$ws.Cells.Item(165,3).Value2 = 1
$ws.Cells.Item(165,5).Value2 = 62
